$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 467, shifting rows 467:551 down to 468:552
$ws.Rows.Item(467).Insert()

# Populate the new row 467 with the new data
$ws.Cells.Item(467, 1).Value = 4
$ws.Cells.Item(467, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(467, 3).Value = "Los Lagos"
$ws.Cells.Item(467, 4).Value = 44711
$ws.Cells.Item(467, 5).Value = 10
$ws.Cells.Item(467, 6).Value = 100112033
$ws.Cells.Item(467, 7).Value = "Lechuga"
$ws.Cells.Item(467, 8).Value = "Escarola"
$ws.Cells.Item(467, 9).Value = "Primera"
$ws.Cells.Item(467, 10).Value = 250
$ws.Cells.Item(467, 11).Value = 12000
$ws.Cells.Item(467, 12).Value = 12000
$ws.Cells.Item(467, 13).Value = 12000
$ws.Cells.Item(467, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(467, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(467, 16).Value = 800
$ws.Cells.Item(467, 17).Value = 15
$ws.Cells.Item(467, 18).Value = "Hortaliza"
